# Apply the "update to published CDA FHIR logical model with patches #241"
# change set to the workbook.
#
# Touches two sheets:
#   - "Metadata": Version / Date / Contact values
#   - "Elements": two ValueSet URLs + Min / Base Min for Order.classCode,
#                 plus the "Binding Value Set" column (Z) grows wider to
#                 fit the new (longer) URL text

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Metadata"
# ---------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

# Version
$meta.Range("B3").Value = "2.0.0-sd-202406-matchbox-patch"

# Date
$meta.Range("B8").Value = "2024-06-19T17:47:42+02:00"

# Contact
$meta.Range("B10").Value = "HL7 International - Structured Documents (http://www.hl7.org/Special/committees/structure, structdog@lists.HL7.org)"

# ---------------------------------------------------------------------
# Sheet "Elements"
# ---------------------------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

# Order.classCode (row 12) - Binding Value Set
$elements.Range("Z12").Value = "http://hl7.org/cda/stds/core/ValueSet/CDAActClass"

# Order.moodCode (row 13) - Binding Value Set
$elements.Range("Z13").Value = "http://hl7.org/cda/stds/core/ValueSet/CDAActMoodIntent"

# Order.classCode Min / Base Min: 1 -> 0 (keep these as *text* cells, like
# the rest of the Min/Max columns in this sheet, not auto-converted numbers)
$elements.Range("F12").NumberFormat = "@"
$elements.Range("F12").Value = "0"

$elements.Range("AG12").NumberFormat = "@"
$elements.Range("AG12").Value = "0"

# Column Z ("Binding Value Set") widens to fit the new, longer URLs
$elements.Range("Z1").EntireColumn.ColumnWidth = 52.8
